# Edit: strip the redundant ".xml" suffixes from the strategy/uncertainty
# add-on file names, and remove the Low-Ambition AFOLU add-on file
# (Strategy_5 / F8:G8) since no add-on file is used for low AFOLU.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# ---- Strategy table (rows 4-9) ----

# Strategy_1 - Decarbonize Power via Renewable Energy
$ws.Range("D4").Value = "Colombia_RPS_High" + $nl + "Colombia_Nuclear_Zero"
$ws.Range("E4").Value = "Strategy_1_High_RPS" + $nl + "Strategy_1_High_Nuclear"
$ws.Range("F4").Value = "Colombia_RPS_Low" + $nl + "Colombia_Nuclear_Normal"
$ws.Range("G4").Value = "Strategy_1_Low_RPS" + $nl + "Strategy_1_Low_Nuclear"

# Strategy_2 - Emphasize Efficiency
$ws.Range("D5").Value = "Colombia_Bldg_ShellApplianceEff_High" + $nl + "Colombia_IndustrialEff_High"
$ws.Range("E5").Value = "Strategy_2_High_BldEE" + $nl + "Strategy_2_High_IndEE"
$ws.Range("F5").Value = "Colombia_Bldg_ShellApplianceEff_Low" + $nl + "Colombia_IndustrialEff_Low"
$ws.Range("G5").Value = "Strategy_2_Low_BldEE" + $nl + "Strategy_2_Low_IndEE"

# Strategy_3 - Electrify Transport
$ws.Range("D6").Value = "transportation_UCD_CORE_RapidEVsw_Colombia"
$ws.Range("E6").Value = "Strategy_3_High_ElecTrans"
$ws.Range("F6").Value = "transportation_UCD_CORE_ModEVsw"
$ws.Range("G6").Value = "Strategy_3_Low_ElecTrans"

# Strategy_4 - Public Transport
$ws.Range("D7").Value = "Colombia_Public_Transport_High"
$ws.Range("E7").Value = "Strategy_4_High_PublicTrans"
$ws.Range("F7").Value = "Colombia_Public_Transport_Normal"
$ws.Range("G7").Value = "Strategy_4_Low_PublicTrans"

# Strategy_5 - AFOLU (Low-ambition add-on file removed - not used)
$ws.Range("D8").Value = "land_constraint_Colombia_10_afforestation"
$ws.Range("E8").Value = "Strategy_5_High_AFOLU"
$ws.Range("F8").Value = $null
$ws.Range("G8").Value = $null

# Strategy_6 - Diet
$ws.Range("D9").Value = "Colombia_Low_Meat"
$ws.Range("E9").Value = "Strategy_6_High_Meat"

# ---- Uncertainty table (rows 13-18) ----

# Uncertaity_1 - Socioeconomics
$ws.Range("D13").Value = "Colombia_GDP_High" + $nl + "Colombia_Population_High"
$ws.Range("E13").Value = "Uncertainty_1_High_GDP" + $nl + "Uncertainty_1_High_Population"
$ws.Range("F13").Value = "Colombia_GDP_Low" + $nl + "Colombia_Population_Low"
$ws.Range("G13").Value = "Uncertainty_1_Low_GDP" + $nl + "Uncertainty_1_Low_Population"

# Uncertaity_2 - EV Costs
$ws.Range("D14").Value = $nl + "transportation_UCD_CORE_RapidEVcost_Colombia_noPubTrninterp"
$ws.Range("E14").Value = "Uncertainty_2_High_EVCost" + $nl
$ws.Range("F14").Value = $nl + "transportation_UCD_CORE_ModEVcost_Colombia_noPubTrninterp"
$ws.Range("G14").Value = "Uncertainty_2_Low_EVCost" + $nl

# Uncertaity_3 - RE Costs (D15 unchanged, E15 loses its ".xml" suffixes)
$ws.Range("E15").Value = "Uncertainty_3_High_RECostSolar Uncertainty_3_High_RECostWind"

# Uncertaity_4 - CCS
$ws.Range("D16").Value = "Global_CCS_Cost_Normal"
$ws.Range("E16").Value = "Uncertainty_4_High_CCSCost"
$ws.Range("F16").Value = "Global_CCS_Cost_High"
$ws.Range("G16").Value = "Uncertainty_4_Low_CCSCost"

# Uncertaity_5 - Climate Change Impacts
$ws.Range("D17").Value = "ag_prodchange_rcp2p6_gfdl_pdssat" + $nl + "hydro_impacts_GFDL-ESM2M_rcp2p6" + $nl + "runoff_impacts_GFDL-ESM2M_rcp2p6"
$ws.Range("E17").Value = "Uncertainty_5_High_Ag" + $nl + "Uncertainty_5_High_Hydro" + $nl + "Uncertainty_5_High_Runnoff"
$ws.Range("F17").Value = "ag_prodchange_rcp2p6_hadgem2_pdssat" + $nl + "hydro_impacts_HadGEM2-ES_rcp2p6" + $nl + "runoff_impacts_HadGEM2-ES_rcp2p6"
$ws.Range("G17").Value = "Uncertainty_5_Low_Ag" + $nl + "Uncertainty_5_Low_Hydro" + $nl + "Uncertainty_5_Low_Runnoff"

# Uncertainty_6 - Trade
$ws.Range("D18").Value = "Global_ag_trade_HOV_CL_25"
$ws.Range("E18").Value = "Uncertainty_6_High_HOV-CL"

# ---- Rich-text formatting for the "NEW High Ambition XML Name" headers ----
# "High" is rendered in a slightly different font (Calibri (Body)); the whole
# "High Ambition XML Name" portion stays bold / dark red.
foreach ($addr in @("E3", "E12")) {
    $cell = $ws.Range($addr)
    $charsHigh = $cell.Characters(5, 4)
    $charsHigh.Font.Name = "Calibri (Body)"
    $charsHigh.Font.Bold = $true
    $charsHigh.Font.Color = 192
    $charsRest = $cell.Characters(9, 18)
    $charsRest.Font.Bold = $true
    $charsRest.Font.Color = 192
}

# ---- Row heights / column widths re-fit after the text shrank ----
$ws.Rows.Item(4).RowHeight = 31.5
$ws.Rows.Item(5).RowHeight = 31.5
$ws.Rows.Item(13).RowHeight = 31.5
$ws.Rows.Item(14).RowHeight = 47.25
$ws.Rows.Item(15).RowHeight = 31.5
$ws.Rows.Item(17).RowHeight = 47.25

$colWidthOffset = 0.8333333333333334
$ws.Columns.Item(3).ColumnWidth = 19.875 - $colWidthOffset
$ws.Columns.Item(4).ColumnWidth = 42.625 - $colWidthOffset
$ws.Columns.Item(5).ColumnWidth = 31.125 - $colWidthOffset
$ws.Columns.Item(6).ColumnWidth = 36.5 - $colWidthOffset
$ws.Columns.Item(7).ColumnWidth = 27.0 - $colWidthOffset

# Column A no longer needs its own explicit width override.
$ws.Columns.Item(1).ColumnWidth = $ws.StandardWidth

# Update the remembered selection to match the last-edited cell.
$ws.Range("H15").Select()
